$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'299.66"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-1.48%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'31.45"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-1.52%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.086"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'-1.89%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.07958"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'6.50%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'2.212"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'-12.75%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'7.745"
$ws.Range("D7").Style = "Normal"
$ws.Range("D8").Value = "'3.862"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.06%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.9165"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-0.18%"
$ws.Range("E9").Style = "Normal"
$ws.Range("E10").Value = "'-0.41%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07385"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-3.13%"
$ws.Range("E11").Style = "Normal"
$ws.Range("D12").Value = "'0.09160"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'11.75%"
$ws.Range("E12").Style = "Normal"
$ws.Range("D13").Value = "'0.03020"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'0.38%"
$ws.Range("E13").Style = "Normal"
$ws.Range("E14").Value = "'0.88%"
$ws.Range("E14").Style = "Normal"
$ws.Range("D15").Value = "'0.001513"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'-0.42%"
$ws.Range("E15").Style = "Normal"
$ws.Range("D16").Value = "'0.005851"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-4.20%"
$ws.Range("E16").Style = "Normal"
$ws.Range("D17").Value = "'3.479"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-0.82%"
$ws.Range("E17").Style = "Normal"
$ws.Range("E18").Value = "'1.67%"
$ws.Range("E18").Style = "Normal"
$ws.Range("E19").Value = "'0.34%"
$ws.Range("E19").Style = "Normal"
$ws.Range("D20").Value = "'0.1336"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.23%"
$ws.Range("E20").Style = "Normal"
$ws.Range("D21").Value = "'4.181"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'-10.27%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'0.1698"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'8.42%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04619"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'0.16%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-1.68%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.004462"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'-1.30%"
$ws.Range("E25").Style = "Normal"
$ws.Range("E26").Value = "'-7.80%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0003394"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'23.75%"
$ws.Range("E27").Style = "Normal"
$ws.Range("D39").Value = "'0.01743"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'-1.58%"
$ws.Range("E39").Style = "Normal"
$ws.Range("D40").Value = "'0.04595"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'1.17%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.006944"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-5.71%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.1356"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-0.37%"
$ws.Range("E42").Style = "Normal"
$ws.Range("E43").Value = "'2.21%"
$ws.Range("E43").Style = "Normal"
$ws.Range("D44").Value = "'0.009553"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'-11.52%"
$ws.Range("E44").Style = "Normal"
$ws.Range("D45").Value = "'0.00006269"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'-3.30%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D46").Value = "'0.00000000749"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'-0.12%"
$ws.Range("E46").Style = "Normal"
$ws.Range("D47").Value = "'0.007972"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'-19.35%"
$ws.Range("E47").Style = "Normal"
$ws.Range("D48").Value = "'1.159"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'41.24%"
$ws.Range("E48").Style = "Normal"
$ws.Range("D49").Value = "'0.00002098"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'-0.12%"
$ws.Range("E49").Style = "Normal"
$ws.Range("D50").Value = "'0.0001998"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'-0.05%"
$ws.Range("E50").Style = "Normal"
